$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.664.03"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "2.200.30"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  +0.09%  "

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "257.45"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.54%  "

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "84.52"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +11.88%  "

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.616"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("E8").Value = "  +0.00%  "

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.597"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "44.85"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +9.52%  "

$ws.Range("E11").Value = "  +0.15%  "

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "7.21"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +4.68%  "

$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("D14").Value = "2.532.43"
$ws.Range("E14").Value = "  -0.39%  "

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "14.39"
$r.Style = "Normal"

$ws.Range("D16").Value = "2.196.96"
$ws.Range("E16").Value = "  -0.75%  "

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.780"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").Value = "43.612.00"
$ws.Range("E18").Value = "  +2.21%  "

$ws.Range("E19").Value = "  +0.19%  "

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "69.82"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "5.91"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("E22").Value = "  +8.09%  "

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "231.69"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.18%  "

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "9.08"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -4.40%  "

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "3.57"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +6.33%  "

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "10.67"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "39.19"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("E30").Value = "  +2.51%  "

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "173.69"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "20.40"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.0860"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +3.67%  "

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "5.35"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +2.82%  "

$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("E36").Value = "  +2.02%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.0360"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +4.37%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "4.49"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +4.71%  "

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "12.43"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("E40").Value = "  +4.45%  "

$ws.Range("E41").Value = "  +0.25%  "

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "62.91"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +4.99%  "

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "5.47"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +3.87%  "

$ws.Range("E44").Value = "  +1.29%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "100.29"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.0980"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "8.30"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("E48").Value = "  +4.83%  "

$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.51"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +9.34%  "

$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.433"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -5.77%  "
